$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(124, 1).Value = 123
$ws.Cells.Item(124, 2).Value = 1
$ws.Cells.Item(124, 3).Value = "2024-06-17 13:15:51"
$ws.Cells.Item(124, 4).Value = 200
$ws.Cells.Item(124, 5).Value = 12

$ws.Cells.Item(125, 1).Value = 124
$ws.Cells.Item(125, 2).Value = 2
$ws.Cells.Item(125, 3).Value = "2024-06-17 13:15:51"
$ws.Cells.Item(125, 4).Value = 200
$ws.Cells.Item(125, 5).Value = 0
